# Technology-cost-database.xlsx edit
#
# The "Database" sheet had several blank rows (5, 7, 9, 11, 14, 16)
# separating technology entries. This edit removes those blank rows
# (shifting the remaining rows up so the data is contiguous: rows 1-11),
# and fixes the Interest Rate ("F") column for the photovoltaic / solar
# collector rows, which previously held a literal "?" placeholder text
# instead of a numeric interest rate -- it is now set to 5% (0.05),
# matching the other technology rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Remove the blank separator rows, starting from the bottom so the
# remaining row numbers referenced below stay valid while deleting.
$ws.Rows(16).Delete() | Out-Null
$ws.Rows(14).Delete() | Out-Null
$ws.Rows(11).Delete() | Out-Null
$ws.Rows(9).Delete()  | Out-Null
$ws.Rows(7).Delete()  | Out-Null
$ws.Rows(5).Delete()  | Out-Null

# After the deletions above, the rows that used to be 12, 13, 15, 17
# (Photovoltaic Panels x2, Solar Collector, PV-thermal) are now rows
# 8, 9, 10, 11. Their Interest Rate ("F") cells held the text "?" --
# replace with the numeric 5% interest rate used elsewhere, formatted
# as a percentage like the other rows.
$interestRateRows = @(8, 9, 10, 11)
foreach ($r in $interestRateRows) {
    $cell = $ws.Range("F" + $r)
    $cell.Value = 0.05
    $cell.NumberFormat = "0%"
}

# Match the saved selection state from the edited workbook.
$ws.Range("D19").Select() | Out-Null
